$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Map of row number -> new value for column G ("Recorded By") cells.
# These updates reorder the comma-separated list of users/emails recorded
# in each row (daily attendance processing normalization).
$updates = @{
    2 = 'backup@backdoor.com, System, system'
    4 = 'backup@backdoor.com, System'
    5 = 'backup@backdoor.com, System'
    8 = 'backup@backdoor.com, System'
    11 = 'dnasr281@gmail.com, System'
    17 = 'dnasr281@gmail.com, System'
    28 = 'backup@backdoor.com, System, system'
    30 = 'backup@backdoor.com, System'
    31 = 'backup@backdoor.com, System'
    34 = 'backup@backdoor.com, System'
    37 = 'dnasr281@gmail.com, System'
    43 = 'dnasr281@gmail.com, System'
    54 = 'backup@backdoor.com, System, system'
    56 = 'backup@backdoor.com, System'
    57 = 'backup@backdoor.com, System'
    60 = 'backup@backdoor.com, System'
    63 = 'dnasr281@gmail.com, System'
    69 = 'dnasr281@gmail.com, System'
    80 = 'backup@backdoor.com, System'
    81 = 'backup@backdoor.com, System'
    82 = 'backup@backdoor.com, System'
    87 = 'dnasr281@gmail.com, admin@admin.com'
    93 = 'dnasr281@gmail.com, System'
    94 = 'dnasr281@gmail.com, System'
    96 = 'dnasr281@gmail.com, System'
    106 = 'backup@backdoor.com, System'
    107 = 'backup@backdoor.com, System'
    108 = 'backup@backdoor.com, System'
    113 = 'dnasr281@gmail.com, admin@admin.com'
    119 = 'dnasr281@gmail.com, System'
    120 = 'dnasr281@gmail.com, System'
    122 = 'dnasr281@gmail.com, System'
    132 = 'backup@backdoor.com, System'
    133 = 'backup@backdoor.com, System'
    134 = 'backup@backdoor.com, System'
    139 = 'dnasr281@gmail.com, admin@admin.com'
    145 = 'dnasr281@gmail.com, System'
    146 = 'dnasr281@gmail.com, System'
    148 = 'dnasr281@gmail.com, System'
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
